$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "30.319.15"
Set-TextValue $ws.Range("E2") "  -3.09%  "
Set-TextValue $ws.Range("D3") "1.935.80"
Set-TextValue $ws.Range("E3") "  -3.19%  "
Set-TextValue $ws.Range("D4") "1.001"
Set-TextValue $ws.Range("D5") "250.31"
Set-TextValue $ws.Range("E5") "  -1.91%  "
Set-TextValue $ws.Range("D6") "0.7219"
Set-TextValue $ws.Range("E6") "  -7.36%  "
Set-TextValue $ws.Range("D8") "0.3308"
Set-TextValue $ws.Range("E8") "  -4.95%  "
Set-TextValue $ws.Range("D9") "27.82"
Set-TextValue $ws.Range("E9") "  -1.41%  "
Set-TextValue $ws.Range("D10") "0.07255"
Set-TextValue $ws.Range("E10") "  +1.89%  "
Set-TextValue $ws.Range("D11") "0.8097"
Set-TextValue $ws.Range("E11") "  -4.25%  "
Set-TextValue $ws.Range("E12") "  -1.48%  "
Set-TextValue $ws.Range("D13") "1.937.23"
Set-TextValue $ws.Range("E13") "  -3.13%  "
Set-TextValue $ws.Range("D14") "5.501"
Set-TextValue $ws.Range("E14") "  -2.92%  "
Set-TextValue $ws.Range("D15") "94.74"
Set-TextValue $ws.Range("E15") "  -6.32%  "
Set-TextValue $ws.Range("D16") "15.10"
Set-TextValue $ws.Range("E16") "  -2.09%  "
Set-TextValue $ws.Range("D17") "30.335.68"
Set-TextValue $ws.Range("D18") "0.000008290"
Set-TextValue $ws.Range("E18") "  +1.59%  "
Set-TextValue $ws.Range("D19") "252.68"
Set-TextValue $ws.Range("E19") "  -7.57%  "
Set-TextValue $ws.Range("D20") "5.890"
Set-TextValue $ws.Range("E20") "  -1.95%  "
Set-TextValue $ws.Range("D21") "2.192.13"
Set-TextValue $ws.Range("E21") "  -2.99%  "
Set-TextValue $ws.Range("E22") "  +0.22%  "
Set-TextValue $ws.Range("D23") "1.001"
Set-TextValue $ws.Range("E23") "  +0.24%  "
Set-TextValue $ws.Range("D24") "6.979"
Set-TextValue $ws.Range("E24") "  -2.16%  "
Set-TextValue $ws.Range("D25") "9.755"
Set-TextValue $ws.Range("E25") "  -3.56%  "
Set-TextValue $ws.Range("D26") "163.78"
Set-TextValue $ws.Range("E26") "  -0.49%  "
Set-TextValue $ws.Range("D27") "2.387"
Set-TextValue $ws.Range("E27") "  -1.01%  "
Set-TextValue $ws.Range("D28") "19.29"
Set-TextValue $ws.Range("E28") "  -3.55%  "
Set-TextValue $ws.Range("D29") "0.1317"
Set-TextValue $ws.Range("E29") "  -7.70%  "
Set-TextValue $ws.Range("D30") "1.568"
Set-TextValue $ws.Range("E30") "  -1.90%  "
Set-TextValue $ws.Range("D31") "1.347"
Set-TextValue $ws.Range("E31") "  -2.12%  "
Set-TextValue $ws.Range("D32") "4.435"
Set-TextValue $ws.Range("E32") "  -4.88%  "
Set-TextValue $ws.Range("D33") "4.180"
Set-TextValue $ws.Range("E33") "  -6.28%  "
Set-TextValue $ws.Range("D34") "0.05202"
Set-TextValue $ws.Range("E34") "  -2.98%  "
Set-TextValue $ws.Range("E35") "  +1.40%  "
Set-TextValue $ws.Range("D36") "0.7507"
Set-TextValue $ws.Range("E36") "  -5.26%  "
Set-TextValue $ws.Range("D37") "2.741"
Set-TextValue $ws.Range("E37") "  -1.14%  "
Set-TextValue $ws.Range("D38") "0.01980"
Set-TextValue $ws.Range("E38") "  -1.53%  "
Set-TextValue $ws.Range("D39") "2.825"
Set-TextValue $ws.Range("E39") "  -3.56%  "
Set-TextValue $ws.Range("D40") "79.38"
Set-TextValue $ws.Range("E40") "  -8.73%  "
Set-TextValue $ws.Range("D41") "6.386"
Set-TextValue $ws.Range("E41") "  -6.17%  "
Set-TextValue $ws.Range("D42") "0.4540"
Set-TextValue $ws.Range("E42") "  -3.43%  "
Set-TextValue $ws.Range("D43") "2.025"
Set-TextValue $ws.Range("E43") "  -5.84%  "
Set-TextValue $ws.Range("D44") "0.8464"
Set-TextValue $ws.Range("E44") "  -1.33%  "
Set-TextValue $ws.Range("D45") "1.0000"
Set-TextValue $ws.Range("E45") "  +0.16%  "
Set-TextValue $ws.Range("D46") "101.80"
Set-TextValue $ws.Range("E46") "  -3.59%  "
Set-TextValue $ws.Range("D47") "9.752"
Set-TextValue $ws.Range("E47") "  -5.42%  "
Set-TextValue $ws.Range("D48") "7.476"
Set-TextValue $ws.Range("E48") "  -4.03%  "
Set-TextValue $ws.Range("D49") "36.86"
Set-TextValue $ws.Range("E49") "  -2.74%  "
Set-TextValue $ws.Range("D50") "0.4190"
Set-TextValue $ws.Range("E50") "  -3.79%  "
Set-TextValue $ws.Range("D51") "0.06034"
Set-TextValue $ws.Range("E51") "  -0.21%  "
